$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.220.49"
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = "'2.214.31"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'296.20"
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").Value = "'87.98"
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("D10").Value = "'52.42"
$ws.Range("E10").Value = '  +7.57%  '
$ws.Range("D11").Value = "'30.94"
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("E13").Value = '  +2.56%  '
$ws.Range("E14").Value = '  -1.41%  '
$ws.Range("D15").Value = "'2.555.53"
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("D16").Value = "'13.87"
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("D17").Value = "'2.213.89"
$ws.Range("E17").Value = '  -2.43%  '
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = "'40.102.66"
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").Value = "'11.34"
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").Value = "'5.78"
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("D23").Value = "'65.75"
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = "'235.61"
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("E27").Value = '  -1.13%  '
$ws.Range("D28").Value = "'23.26"
$ws.Range("E28").Value = '  +2.14%  '
$ws.Range("E29").Value = '  +0.96%  '
$ws.Range("E30").Value = '  -4.90%  '
$ws.Range("D31").Value = "'156.42"
$ws.Range("D32").Value = "'32.21"
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("E35").Value = '  +3.32%  '
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("D37").Value = "'2.33"
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("D40").Value = "'1.74"
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("D41").Value = "'15.63"
$ws.Range("E41").Value = '  -1.36%  '
$ws.Range("D42").Value = "'3.83"
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("D43").Value = "'2.070.52"
$ws.Range("E43").Value = '  -2.89%  '
$ws.Range("D44").Value = "'19.33"
$ws.Range("E44").Value = '  +5.39%  '
$ws.Range("E45").Value = '  +0.61%  '
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("D47").Value = "'2.83"
$ws.Range("E47").Value = '  +6.05%  '
$ws.Range("D48").Value = "'1.91"
$ws.Range("E48").Value = '  -11.60%  '
$ws.Range("D49").Value = "'2.428.40"
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("E50").Value = '  +1.95%  '
$ws.Range("E51").Value = '  +0.22%  '
